$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.762.31'
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.853.02'
$ws.Range("E4").Value = '  -1.43%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.00'
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("E6").Value = '  -1.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4303'
$ws.Range("E7").Value = '  -1.65%  '
$ws.Range("E8").Value = '  -0.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07342'
$ws.Range("E9").Value = '  -0.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8777'
$ws.Range("E10").Value = '  -0.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.58'
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.857.13'
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.760'
$ws.Range("E13").Value = '  +0.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.445'
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07118'
$ws.Range("E15").Value = '  +0.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.96'
$ws.Range("E16").Value = '  +4.62%  '
$ws.Range("E17").Value = '  -1.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009012'
$ws.Range("E18").Value = '  -0.79%  '
$ws.Range("E19").Value = '  -1.35%  '
$ws.Range("E20").Value = '  +0.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.790.18'
$ws.Range("E21").Value = '  +0.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.224'
$ws.Range("E22").Value = '  -1.11%  '
$ws.Range("E23").Value = '  -1.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.083.61'
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.985'
$ws.Range("E25").Value = '  -1.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.46'
$ws.Range("E26").Value = '  -1.11%  '
$ws.Range("E27").Value = '  -0.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.179'
$ws.Range("E28").Value = '  +9.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.364'
$ws.Range("E29").Value = '  +0.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.98'
$ws.Range("E30").Value = '  +1.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08942'
$ws.Range("E31").Value = '  -0.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.230'
$ws.Range("E32").Value = '  +1.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7780'
$ws.Range("E33").Value = '  +0.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.547'
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.932'
$ws.Range("E35").Value = '  -1.98%  '
$ws.Range("E36").Value = '  -1.31%  '
$ws.Range("E37").Value = '  -0.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01985'
$ws.Range("E38").Value = '  +0.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05348'
$ws.Range("E39").Value = '  +1.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.907'
$ws.Range("E40").Value = '  +1.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.166'
$ws.Range("E41").Value = '  +4.79%  '
$ws.Range("E42").Value = '  +1.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5136'
$ws.Range("E43").Value = '  -0.87%  '
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.78'
$ws.Range("E45").Value = '  +0.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '107.59'
$ws.Range("E46").Value = '  -2.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4767'
$ws.Range("E47").Value = '  +1.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06477'
$ws.Range("E48").Value = '  -1.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.693'
$ws.Range("E49").Value = '  -0.48%  '
$ws.Range("E50").Value = '  -1.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.844'
